$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3: "hackathon" -> "Hackathon"
$ws.Range("C3").Value = "Hackathon"

# C4: "life" -> "Debate"
$ws.Range("C4").Value = "Debate"

# D4: text "Participation" -> numeric 3
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Value = 3

# Match final selection/active cell state left after editing
$null = $ws.Range("C12").Select()
